$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 corresponds to ICSA_thou / Initial Jobless Claims - auto-updated data refresh
$ws.Range("E9").Value = 206000
$ws.Range("G9").Value = 363777.7777777778
$ws.Range("H9").Value = -9000
$ws.Range("I9").Value = -0.04186046511627907
